$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New phone numbers for rows 4 and 7
$ws.Range("C4").Value = 89036054939
$ws.Range("C7").Value = 81233214323

# New email for B5 (cell already carried the hyperlink style s="1" but was
# empty) - set the value, add the hyperlink, then restore the existing
# hyperlink cell-style formatting (Excel's Hyperlinks.Add otherwise stamps a
# brand-new style index instead of reusing the workbook's existing one).
$ws.Range("B5").Value = "awdawd@test.com"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:awdawd@test.com")
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# New email for B7 (brand-new cell in a row that previously only had A7)
$ws.Range("B7").Value = "dawdaw@tesd.tv"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:dawdaw@tesd.tv")
$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update selection to C7
$ws.Range("C7").Select()
